# TC05_Search_product_in_Catalog.xlsx -- "Changes for New UI Prod"
#
# Sheet "Testdata" gets two brand-new rows (EleType1/JSElement,
# EleType2/JSElement), the old 203kdd / 203KDD_RADIAL BALL BEARING sample
# values are replaced by a numeric dimension (200) and a boolean flag, and
# the "validSearchText" data-descriptor is renamed to "Pagination".
#
# Sheet "TC05_Search_product_in_Catalog" swaps the old
# VERIFY_TEXT_PRESENT / ValidSearchHeader / validSearchText row for a new
# VERIFY_WEBELEMENT_PRESENT / ValidSearchPagination / Pagination row, the
# ValidSearchHeader object on the row above becomes SearchDimensions, and
# the trailing duplicate VERIFY_WEBELEMENT_PRESENT row (row 8) is removed
# entirely.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TC05_Search_product_in_Catalog")
$ws2 = $wb.Worksheets.Item("Testdata")

# --- Testdata: add the two new data-object rows -----------------------
$ws2.Range("A7").Value = "EleType1"
$ws2.Range("B7").Value = "JSElement"
$ws2.Range("A8").Value = "EleType2"
$ws2.Range("B8").Value = "JSElement"

# --- TC05_Search_product_in_Catalog: relabel the search-dimensions row -
$ws1.Range("C6").Value = "SearchDimensions"

# --- TC05_Search_product_in_Catalog: replace the old text-present check
#     (row 7) with the new pagination web-element check ----------------
$ws1.Range("B7").Value = "VERIFY_WEBELEMENT_PRESENT"
$ws1.Range("C7").Value = "ValidSearchPagination"
$ws1.Range("E7").Value = "Pagination"

# --- Testdata: rename validSearchText -> Pagination, and swap the old
#     string sample values for a numeric width and a boolean flag -------
$ws2.Range("B3").Value = 200
$ws2.Range("A5").Value = "Pagination"
$ws2.Range("B5").Value = $true

# --- TC05_Search_product_in_Catalog: row 8 (duplicate
#     VERIFY_WEBELEMENT_PRESENT / ValidSeachImg check) is gone -----------
[void]$ws1.Rows.Item(8).Delete()

# Row 4 picks up a slightly tighter custom row height in the new layout.
$ws1.Rows.Item(4).RowHeight = 14.25

# --- selection state: Testdata keeps B5 selected, the test-case sheet
#     is the active tab with A3:XFD7 selected -----------------------
[void]$ws2.Range("B5").Select()
[void]$ws1.Range("A3:XFD7").Select()
